$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"

# --- Narrow the status columns (was ~17.22 chars, now ~13.41 chars) ---
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5
